# Update the "想去人数" (F column) counts for the sheets "展览" (1) and "全部类型" (4).
# Row -> [old_sheet1_value, new_sheet1_value, old_sheet4_value, new_sheet4_value]
$wb = $excel.ActiveWorkbook

$updates = @(
    @{Row=3;  S1New=1446; S4New=1446},
    @{Row=4;  S1New=175;  S4New=175},
    @{Row=8;  S1New=21;   S4New=21},
    @{Row=10; S1New=142;  S4New=142},
    @{Row=11; S1New=5;    S4New=5},
    @{Row=12; S1New=4784; S4New=4784},
    @{Row=14; S1New=7069; S4New=7069},
    @{Row=18; S1New=583;  S4New=583},
    @{Row=20; S1New=9;    S4New=9},
    @{Row=21; S1New=4190; S4New=4190},
    @{Row=22; S1New=1320; S4New=1320},
    @{Row=24; S1New=78;   S4New=78},
    @{Row=25; S1New=2770; S4New=2770},
    @{Row=26; S1New=576;  S4New=576},
    @{Row=27; S1New=557;  S4New=557},
    @{Row=28; S1New=178;  S4New=178},
    @{Row=29; S1New=396;  S4New=396},
    @{Row=31; S1New=418;  S4New=418},
    @{Row=33; S1New=56;   S4New=56},
    @{Row=35; S1New=1069; S4New=1069},
    @{Row=36; S1New=75;   S4New=75},
    @{Row=37; S1New=923;  S4New=923},
    @{Row=39; S1New=559;  S4New=559},
    @{Row=40; S1New=12;   S4New=12},
    @{Row=44; S1New=94;   S4New=94},
    @{Row=45; S1New=852;  S4New=853},
    @{Row=46; S1New=660;  S4New=660},
    @{Row=47; S1New=26;   S4New=26}
)

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.S1New
    $ws4.Cells.Item($u.Row, 6).Value = $u.S4New
}
